$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 8
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(6, 6).Value = 4
$ws.Cells.Item(7, 6).Value = 4
$ws.Cells.Item(8, 6).Value = 4
$ws.Cells.Item(10, 6).Value = 5
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 5
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 7
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(19, 5).Value = 18
$ws.Cells.Item(19, 6).Value = 3
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 8
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(22, 6).Value = 7
$ws.Cells.Item(23, 6).Value = 1
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 6).Value = 8
$ws.Cells.Item(25, 5).Value = 2
$ws.Cells.Item(25, 6).Value = 11
$ws.Cells.Item(26, 5).Value = 1
$ws.Cells.Item(26, 6).Value = 5
$ws.Cells.Item(27, 5).ClearContents()
$ws.Cells.Item(27, 6).ClearContents()
$ws.Cells.Item(28, 6).Value = 2
$ws.Cells.Item(29, 5).Value = 9
$ws.Cells.Item(30, 5).Value = 2
$ws.Cells.Item(30, 6).Value = 17
$ws.Cells.Item(31, 5).Value = 6
$ws.Cells.Item(31, 6).Value = 23
$ws.Cells.Item(32, 5).Value = 1
$ws.Cells.Item(32, 6).Value = 15
$ws.Cells.Item(34, 5).Value = 6
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 14
$ws.Cells.Item(36, 5).Value = 3
$ws.Cells.Item(36, 6).Value = 10
$ws.Cells.Item(37, 5).Value = 2
$ws.Cells.Item(37, 6).Value = 15
$ws.Cells.Item(38, 5).Value = 1
$ws.Cells.Item(38, 6).Value = 12
$ws.Cells.Item(39, 5).Value = 6
$ws.Cells.Item(40, 5).Value = 4
$ws.Cells.Item(40, 6).Value = 30
$ws.Cells.Item(41, 6).Value = 8
$ws.Cells.Item(43, 5).Value = 2
$ws.Cells.Item(44, 6).Value = 2
$ws.Cells.Item(45, 5).ClearContents()
$ws.Cells.Item(45, 6).ClearContents()
$ws.Cells.Item(46, 6).Value = 1
$ws.Cells.Item(47, 5).Value = 4
$ws.Cells.Item(47, 6).Value = 13
$ws.Cells.Item(48, 5).Value = 0
$ws.Cells.Item(48, 6).Value = 4
$ws.Cells.Item(49, 5).Value = 2
$ws.Cells.Item(49, 6).Value = 7
$ws.Cells.Item(51, 6).Value = 7
$ws.Cells.Item(52, 5).Value = 1
$ws.Cells.Item(52, 6).Value = 12
$ws.Cells.Item(53, 6).Value = 3
$ws.Cells.Item(54, 5).Value = 16
$ws.Cells.Item(54, 6).Value = 2
$ws.Cells.Item(55, 5).Value = 0
$ws.Cells.Item(55, 6).Value = 9
$ws.Cells.Item(57, 5).Value = 3
$ws.Cells.Item(59, 5).Value = 0
$ws.Cells.Item(60, 5).Value = 0
$ws.Cells.Item(60, 6).Value = 6
$ws.Cells.Item(61, 5).Value = 0
$ws.Cells.Item(61, 6).Value = 11
$ws.Cells.Item(62, 5).Value = 1
$ws.Cells.Item(62, 6).Value = 5
$ws.Cells.Item(63, 5).Value = 2
$ws.Cells.Item(64, 5).Value = 1
$ws.Cells.Item(64, 6).Value = 9
$ws.Cells.Item(65, 5).Value = 2
$ws.Cells.Item(65, 6).Value = 7
$ws.Cells.Item(66, 5).Value = 1
$ws.Cells.Item(67, 5).Value = 24
$ws.Cells.Item(67, 6).Value = 8
$ws.Cells.Item(68, 5).Value = 1
$ws.Cells.Item(68, 6).Value = 9
$ws.Cells.Item(69, 5).Value = 3
$ws.Cells.Item(69, 6).Value = 10
$ws.Cells.Item(70, 6).Value = 2
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(71, 6).Value = 7
$ws.Cells.Item(72, 5).Value = 1
$ws.Cells.Item(72, 6).Value = 11
$ws.Cells.Item(73, 5).Value = 1
$ws.Cells.Item(73, 6).Value = 2
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(75, 6).Value = 3
$ws.Cells.Item(76, 5).Value = 1
$ws.Cells.Item(76, 6).Value = 12
$ws.Cells.Item(77, 5).Value = 4
$ws.Cells.Item(77, 6).Value = 4
$ws.Cells.Item(78, 5).Value = 0
$ws.Cells.Item(78, 6).Value = 5
$ws.Cells.Item(79, 5).Value = 1
$ws.Cells.Item(79, 6).Value = 8
$ws.Cells.Item(80, 5).Value = 0
$ws.Cells.Item(81, 5).Value = 0
$ws.Cells.Item(81, 6).Value = 5
$ws.Cells.Item(82, 6).Value = 7
$ws.Cells.Item(83, 5).Value = 0
$ws.Cells.Item(83, 6).Value = 3
$ws.Cells.Item(84, 5).Value = 7
$ws.Cells.Item(84, 6).Value = 22
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 4
$ws.Cells.Item(86, 5).Value = 2
$ws.Cells.Item(86, 6).Value = 5
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 5).Value = 2
$ws.Cells.Item(89, 5).Value = 3
$ws.Cells.Item(89, 6).Value = 8
$ws.Cells.Item(90, 5).Value = 0
$ws.Cells.Item(90, 6).Value = 2
$ws.Cells.Item(91, 5).Value = 0
$ws.Cells.Item(91, 6).Value = 4
